# Update C2 timestamp, then append new error-log rows (3-11) for the
# additional "pastas" (folders) that failed to process, per the commit:
# "corrigindo o arquivo - processar pastas / resolvendo erro de
# atualização da pasta geral".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2..11: Pasta, Erro, DataHora
$rows = @(
    @{ Row = 2;  Pasta = "BASE";     Hora = "24/07/2025 13:29:39" },
    @{ Row = 3;  Pasta = "CONDONAL"; Hora = "24/07/2025 13:29:39" },
    @{ Row = 4;  Pasta = "FOLK";     Hora = "24/07/2025 13:29:39" },
    @{ Row = 5;  Pasta = "GA";       Hora = "24/07/2025 13:29:39" },
    @{ Row = 6;  Pasta = "GESTART";  Hora = "24/07/2025 13:29:53" },
    @{ Row = 7;  Pasta = "HISEG";    Hora = "24/07/2025 13:29:53" },
    @{ Row = 8;  Pasta = "PRIMEE";   Hora = "24/07/2025 13:29:53" },
    @{ Row = 9;  Pasta = "SINGULAR"; Hora = "24/07/2025 13:29:53" },
    @{ Row = 10; Pasta = "UNICA";    Hora = "24/07/2025 13:29:53" },
    @{ Row = 11; Pasta = "VIGON";    Hora = "24/07/2025 13:29:53" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $pasta = $item.Pasta
    $erro = "Destinatário não encontrado para '" + $pasta + "'."

    $ws.Cells.Item($r, 1).Value = $pasta
    $ws.Cells.Item($r, 2).Value = $erro
    $ws.Cells.Item($r, 3).Value = $item.Hora
}
